# Auto-generated Word COM-interop script implementing the diff
$d = $word.ActiveDocument

# Paragraph 6 (1/28/14): remove _GoBack bookmark
$p = $d.Paragraphs.Item(6)
$r = $d.Range($p.Range.Start, $p.Range.End)
$xml6 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1/28/14</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml6)

# Paragraph 7 (Fortran day 1): merge runs, drop proofErr
$p = $d.Paragraphs.Item(7)
$r = $d.Range($p.Range.Start, $p.Range.End)
$xml7 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">I spent 3 hours today working on the Fortran program.  At first for the first hour I spent just trying to compile the language as well as look at the basic data types.  After I set up my git I tested how to manipulate strings and other types.  I realized that strings are terrible in Fortran.  They are just an array of characters which is not as easy as I would like to work with.  Instead of index’s you need to substring them to get a character to work with.  It took me about 45 minutes to work on the algorithm and such for the functions, but I struggled to insert the characters one grabbed from the string back into the string (Array of Characters).  After a lot of research I figured out a pretty simple way to do it, which was frustrating because of how simple it was and how hard it was to find the answer.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml7)

# Paragraph 8 (Overall I liked Fortran...): merge first part, keep rest
$p = $d.Paragraphs.Item(8)
$r = $d.Range($p.Range.Start, $p.Range.End)
$xml8 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Overall I liked Fortran for this assignment.  I would imagine from the small parts I used it, it would be hard to write a scalable program in this language but I could see how since I feel you have much control over the language it could be easy to do certain things.  It reminded me of sML I do not know if that was because of how I wrote it in subroutines and that is sort of how functions are in sML, but I did not hate Fortran once I began to learn how to use it.  </w:t></w:r><w:r><w:t xml:space="preserve">The readability of it is pretty simple, especially to a programmer.  To a </w:t></w:r><w:r><w:t>non-programmer</w:t></w:r><w:r><w:t xml:space="preserve"> certain </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">things like substring might look weird.  I enjoyed the start and end instead of {} because it makes it easier to read in my opinion.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml8)

# Paragraph 10 (cobol): merge runs, drop proofErr
$p = $d.Paragraphs.Item(10)
$r = $d.Range($p.Range.Start, $p.Range.End)
$xml10 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">I spent a few hours today trying to research cobol, with very little luck.  I had an extremely difficult time trying to find good sources of code to examine to help myself with trying to declare variables as well as understand what limitations the language has.  I eventually got to the point of declaring variables and making a loop.  I grew to understand what exactly I was doing and how I could do things.  I wrote the loops needed to encode, decode, and solve the string. The next issue I faced was making it its own function or subroutine. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml10)

# Paragraph 12 (Today I spent about an hour... cobol day2): merge runs, drop proofErr
$p = $d.Paragraphs.Item(12)
$r = $d.Range($p.Range.Start, $p.Range.End)
$xml12 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">Today I spent about an hour trying to research how to tackle the subroutines needed to make each their own function.  After a lot of time, I found a good source of code I could read to help myself with this issue.  I then after a few minutes separated the code into three routines, and fixed up my code.  For some reason comments give me errors with the online compiler.  Also I realized I had a bug in the way I declared one of my variables.  I only had PIC 99 so the digits in that integer did not go into the hundreds column so I needed to add another 9 which limited the way the alphabet looped around since when you add 26 to some letters it gets above 100 which caused an issue because before this it went back to 0. Overall I did not mind the language once I got around to seeing how it worked.  If there was better documentation or resources online this part would not have been so bad.  The language is pretty straight forward once you read it.   I like certain aspects like explicitly stating end as well as display instead of like System.out.println in java.  When to use periods and when not to is very confusing as well as I do not like the </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">MOVE keyword instead of using =.  I am used to the = vs == so it does not bother me but I have begun to look into pascal and I like the := operator to setting values to a variable.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml12)

# Paragraph 13 (I continued to work today... pascal): merge runs, drop proofErr
$p = $d.Paragraphs.Item(13)
$r = $d.Range($p.Range.Start, $p.Range.End)
$xml13 = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">I continued to work today and finished the pascal.  It was by far the easiest of the three I have done so far.  It reminds me a lot sML.  The var to declare things before the main part of the method is like sML and the way it feels when writing it feels similar as well.  There was a lot of good resources online to help which made it very easy.  It maybe took me over an hour to write, since I know the algorithm so well now.  All it really takes is the ability declare data types and variables,  loop, substring a string, change the character to ASCII then change it back and re-adding it to a string.  Once I figure out how to do those tasks the rest of the program is extremely easy.  </w:t></w:r><w:r><w:t>This was not a bad language to learn, I really enjoyed it.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml13)

# Append two new paragraphs (2/6/14 date + scala entry with new _GoBack bookmark)
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$rEnd = $d.Range($pLast.Range.End, $pLast.Range.End)
$xmlNew = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>2/6/14</w:t></w:r></w:p><w:p><w:r><w:t>Today I did scala procedural.  It was not hard at all, it had a lot of documentation online as well as it was easy to compile.  It is very similar to java, just a little different syntax wise.  I enjoyed how easy it was and how there was type casting.  It caused a little but of an issue going from chars to strings, but after figuring that out it was fairly simple.  This was my favorite language that I have done so far for this assignment.  I think it</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> was because it reminded me of java, which is my favorite language.  </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rEnd.InsertXML($xmlNew)

Write-Host "Final paragraph count:" $d.Paragraphs.Count
